# Insert a new weekly price record for Papaya (Vega Modelo de Temuco) as
# row 14, pushing the existing rows 14-68 down to 15-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 (shifts old row 14.. down by one).
$ws.Rows.Item(14).Insert()

# Populate the new row with the new observation.
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44620
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100108
$ws.Cells.Item(14, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(14, 9).Value = 100108004
$ws.Cells.Item(14, 10).Value = "Papaya"
$ws.Cells.Item(14, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 35
$ws.Cells.Item(14, 14).Value = 24000
$ws.Cells.Item(14, 15).Value = 24000
$ws.Cells.Item(14, 16).Value = 24000
$ws.Cells.Item(14, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 19).Value = 2400
$ws.Cells.Item(14, 20).Value = 10
